# Auto-generated Excel COM-interop script
# Applies the "Horarios actualizados Linea 141 - 297" data refresh:
#  - New scrape timestamp (Ultima actualizacion: 20:46:10) on all 3 sheets
#  - 9 new schedule rows appended to sheet "LP1912" (now sorted by Hora_Llegada)
#  - 1 new schedule row appended to sheet "LP1912-215"
#  - "Total filas" counters updated on sheet "LP1912" (533) and "LP1912-215" (56)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Sheet: LP1912 ---
$ws1.Cells.Item(2, 1).Value = 'Última actualización: 20:46:10'
$ws1.Cells.Item(3, 1).Value = 'Total filas: 533'
$ws1.Cells.Item(121, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(122, 1).Value = '09:38:09'
$ws1.Cells.Item(122, 3).Value = '14_ABASTO'
$ws1.Cells.Item(122, 4).Value = 3
$ws1.Cells.Item(123, 1).Value = '08:21:50'
$ws1.Cells.Item(123, 3).Value = '215C_EL PATO'
$ws1.Cells.Item(123, 4).Value = 80
$ws1.Cells.Item(204, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(205, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(253, 1).Value = '12:37:14'
$ws1.Cells.Item(253, 3).Value = '16_P MOR-167 Y 521'
$ws1.Cells.Item(253, 4).Value = 80
$ws1.Cells.Item(254, 1).Value = '13:53:08'
$ws1.Cells.Item(254, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(254, 4).Value = 4
$ws1.Cells.Item(288, 3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(289, 3).Value = '15_ABASTO'
$ws1.Cells.Item(323, 1).Value = '15:31:33'
$ws1.Cells.Item(323, 3).Value = '215C_EL PATO'
$ws1.Cells.Item(323, 4).Value = 49
$ws1.Cells.Item(324, 1).Value = '14:46:52'
$ws1.Cells.Item(324, 3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(324, 4).Value = 94
$ws1.Cells.Item(339, 3).Value = '225_GOMEZ'
$ws1.Cells.Item(340, 3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(373, 1).Value = '16:43:37'
$ws1.Cells.Item(373, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(373, 4).Value = 51
$ws1.Cells.Item(374, 1).Value = '16:53:01'
$ws1.Cells.Item(374, 3).Value = '10_OLMOS'
$ws1.Cells.Item(374, 4).Value = 41
$ws1.Cells.Item(405, 1).Value = '17:41:19'
$ws1.Cells.Item(405, 3).Value = '10_OLMOS'
$ws1.Cells.Item(405, 4).Value = 35
$ws1.Cells.Item(406, 1).Value = '17:14:55'
$ws1.Cells.Item(406, 3).Value = '15_ABASTO'
$ws1.Cells.Item(406, 4).Value = 62
$ws1.Cells.Item(452, 1).Value = '17:41:19'
$ws1.Cells.Item(452, 3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(452, 4).Value = 100
$ws1.Cells.Item(453, 1).Value = '18:37:25'
$ws1.Cells.Item(453, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(453, 4).Value = 44
$ws1.Cells.Item(454, 3).Value = '14_ABASTO'
$ws1.Cells.Item(470, 3).Value = '81_EL PELIGRO'
$ws1.Cells.Item(471, 1).Value = '18:37:25'
$ws1.Cells.Item(471, 3).Value = '11X44_ETCHEVERRY'
$ws1.Cells.Item(471, 4).Value = 74
$ws1.Cells.Item(472, 1).Value = '18:17:05'
$ws1.Cells.Item(472, 3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(472, 4).Value = 94
$ws1.Cells.Item(479, 1).Value = '19:56:21'
$ws1.Cells.Item(479, 3).Value = '14_ABASTO'
$ws1.Cells.Item(479, 4).Value = 4
$ws1.Cells.Item(480, 1).Value = '18:37:25'
$ws1.Cells.Item(480, 3).Value = '17_ROMERO'
$ws1.Cells.Item(480, 4).Value = 83
$ws1.Cells.Item(488, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(489, 3).Value = '15_ABASTO'
$ws1.Cells.Item(490, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(500, 3).Value = '15_ABASTO'
$ws1.Cells.Item(501, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(507, 1).Value = '20:46:10'
$ws1.Cells.Item(507, 2).Value = '20:46'
$ws1.Cells.Item(507, 3).Value = '10_OLMOS'
$ws1.Cells.Item(507, 4).Value = 0
$ws1.Cells.Item(508, 1).Value = '20:46:10'
$ws1.Cells.Item(508, 2).Value = '20:47'
$ws1.Cells.Item(508, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(508, 4).Value = 1
$ws1.Cells.Item(509, 1).Value = '18:58:44'
$ws1.Cells.Item(509, 2).Value = '20:48'
$ws1.Cells.Item(509, 3).Value = '14X44_ABASTO'
$ws1.Cells.Item(509, 4).Value = 110
$ws1.Cells.Item(510, 1).Value = '19:56:21'
$ws1.Cells.Item(510, 2).Value = '20:52'
$ws1.Cells.Item(510, 4).Value = 56
$ws1.Cells.Item(511, 1).Value = '19:42:02'
$ws1.Cells.Item(511, 2).Value = '20:52'
$ws1.Cells.Item(511, 3).Value = '15_ABASTO'
$ws1.Cells.Item(511, 4).Value = 70
$ws1.Cells.Item(512, 2).Value = '20:53'
$ws1.Cells.Item(512, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(512, 4).Value = 71
$ws1.Cells.Item(513, 1).Value = '20:46:10'
$ws1.Cells.Item(513, 2).Value = '20:54'
$ws1.Cells.Item(513, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(513, 4).Value = 8
$ws1.Cells.Item(514, 1).Value = '18:58:44'
$ws1.Cells.Item(514, 2).Value = '20:56'
$ws1.Cells.Item(514, 3).Value = '10_OLMOS'
$ws1.Cells.Item(514, 4).Value = 118
$ws1.Cells.Item(515, 2).Value = '20:57'
$ws1.Cells.Item(515, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(515, 4).Value = 75
$ws1.Cells.Item(516, 1).Value = '19:56:21'
$ws1.Cells.Item(516, 2).Value = '21:00'
$ws1.Cells.Item(516, 3).Value = '215B_EL PATO'
$ws1.Cells.Item(516, 4).Value = 64
$ws1.Cells.Item(517, 2).Value = '21:01'
$ws1.Cells.Item(517, 3).Value = '215B_EL PATO'
$ws1.Cells.Item(517, 4).Value = 79
$ws1.Cells.Item(518, 2).Value = '21:04'
$ws1.Cells.Item(518, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(518, 4).Value = 82
$ws1.Cells.Item(519, 1).Value = '20:15:00'
$ws1.Cells.Item(519, 2).Value = '21:16'
$ws1.Cells.Item(519, 4).Value = 61
$ws1.Cells.Item(520, 1).Value = '19:42:02'
$ws1.Cells.Item(520, 2).Value = '21:21'
$ws1.Cells.Item(520, 3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(520, 4).Value = 99
$ws1.Cells.Item(521, 2).Value = '21:23'
$ws1.Cells.Item(521, 3).Value = '10_OLMOS'
$ws1.Cells.Item(521, 4).Value = 101
$ws1.Cells.Item(522, 1).Value = '20:46:10'
$ws1.Cells.Item(522, 2).Value = '21:23'
$ws1.Cells.Item(522, 3).Value = '15_ABASTO'
$ws1.Cells.Item(522, 4).Value = 37
$ws1.Cells.Item(523, 1).Value = '20:31:05'
$ws1.Cells.Item(523, 2).Value = '21:31'
$ws1.Cells.Item(523, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(523, 4).Value = 60
$ws1.Cells.Item(524, 1).Value = '20:46:10'
$ws1.Cells.Item(524, 2).Value = '21:32'
$ws1.Cells.Item(524, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(524, 4).Value = 46
$ws1.Cells.Item(525, 2).Value = '21:33'
$ws1.Cells.Item(525, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(525, 4).Value = 62
$ws1.Cells.Item(526, 1).Value = '20:46:10'
$ws1.Cells.Item(526, 2).Value = '21:34'
$ws1.Cells.Item(526, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(526, 4).Value = 48
$ws1.Cells.Item(527, 1).Value = '19:42:02'
$ws1.Cells.Item(527, 2).Value = '21:38'
$ws1.Cells.Item(527, 3).Value = '17_ROMERO'
$ws1.Cells.Item(527, 4).Value = 116
$ws1.Cells.Item(528, 1).Value = '19:42:02'
$ws1.Cells.Item(528, 2).Value = '21:38'
$ws1.Cells.Item(528, 3).Value = '14_ABASTO'
$ws1.Cells.Item(528, 4).Value = 116
$ws1.Cells.Item(529, 1).Value = '20:15:00'
$ws1.Cells.Item(529, 2).Value = '21:43'
$ws1.Cells.Item(529, 3).Value = '17_ROMERO'
$ws1.Cells.Item(529, 4).Value = 88
$ws1.Cells.Item(530, 1).Value = '19:56:21'
$ws1.Cells.Item(530, 2).Value = '21:47'
$ws1.Cells.Item(530, 3).Value = '215A_EL PATO'
$ws1.Cells.Item(530, 4).Value = 111
$ws1.Cells.Item(530, 5).Value = 'LP1912'
$ws1.Cells.Item(531, 1).Value = '20:31:05'
$ws1.Cells.Item(531, 2).Value = '21:58'
$ws1.Cells.Item(531, 3).Value = '17_ROMERO'
$ws1.Cells.Item(531, 4).Value = 87
$ws1.Cells.Item(531, 5).Value = 'LP1912'
$ws1.Cells.Item(532, 1).Value = '20:15:00'
$ws1.Cells.Item(532, 2).Value = '22:08'
$ws1.Cells.Item(532, 3).Value = '17_ROMERO'
$ws1.Cells.Item(532, 4).Value = 113
$ws1.Cells.Item(532, 5).Value = 'LP1912'
$ws1.Cells.Item(533, 1).Value = '20:31:05'
$ws1.Cells.Item(533, 2).Value = '22:08'
$ws1.Cells.Item(533, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(533, 4).Value = 97
$ws1.Cells.Item(533, 5).Value = 'LP1912'
$ws1.Cells.Item(534, 1).Value = '20:46:10'
$ws1.Cells.Item(534, 2).Value = '22:14'
$ws1.Cells.Item(534, 3).Value = '17_ROMERO'
$ws1.Cells.Item(534, 4).Value = 88
$ws1.Cells.Item(534, 5).Value = 'LP1912'
$ws1.Cells.Item(535, 1).Value = '20:46:10'
$ws1.Cells.Item(535, 2).Value = '22:21'
$ws1.Cells.Item(535, 3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(535, 4).Value = 95
$ws1.Cells.Item(535, 5).Value = 'LP1912'
$ws1.Cells.Item(536, 1).Value = '20:31:05'
$ws1.Cells.Item(536, 2).Value = '22:23'
$ws1.Cells.Item(536, 3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(536, 4).Value = 112
$ws1.Cells.Item(536, 5).Value = 'LP1912'
$ws1.Cells.Item(537, 1).Value = '20:31:05'
$ws1.Cells.Item(537, 2).Value = '22:28'
$ws1.Cells.Item(537, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(537, 4).Value = 117
$ws1.Cells.Item(537, 5).Value = 'LP1912'
$ws1.Cells.Item(538, 1).Value = '20:46:10'
$ws1.Cells.Item(538, 2).Value = '22:39'
$ws1.Cells.Item(538, 3).Value = '215A_EL PATO'
$ws1.Cells.Item(538, 4).Value = 113
$ws1.Cells.Item(538, 5).Value = 'LP1912'
# --- Sheet: LP1912-215 ---
$ws2.Cells.Item(2, 1).Value = 'Última actualización: 20:46:10'
$ws2.Cells.Item(3, 1).Value = 'Total filas: 56'
$ws2.Cells.Item(61, 1).Value = '20:46:10'
$ws2.Cells.Item(61, 2).Value = '22:39'
$ws2.Cells.Item(61, 3).Value = '215A_EL PATO'
$ws2.Cells.Item(61, 4).Value = 113
$ws2.Cells.Item(61, 5).Value = 'LP1912'
# --- Sheet: 6203-6173 ---
$ws3.Cells.Item(2, 1).Value = 'Última actualización: 20:46:10'
